$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 291; this shifts existing rows 291:331 down to 292:332
$ws.Rows.Item(291).Insert()

# Populate the newly inserted row 291 with the new record's data.
# Columns A,B,C,E,F,G,H,I,N,Q,R follow the same pattern as the surrounding rows.
$ws.Cells.Item(291, 1).Value = 9
$ws.Cells.Item(291, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(291, 3).Value = "Metropolitana"
$ws.Cells.Item(291, 4).Value = 45077
$ws.Cells.Item(291, 5).Value = 13
$ws.Cells.Item(291, 6).Value = 100112026
$ws.Cells.Item(291, 7).Value = "Haba"
$ws.Cells.Item(291, 8).Value = "Sin especificar"
$ws.Cells.Item(291, 9).Value = "Primera"
$ws.Cells.Item(291, 10).Value = 70
$ws.Cells.Item(291, 11).Value = 16000
$ws.Cells.Item(291, 12).Value = 18000
$ws.Cells.Item(291, 13).Value = 17000
$ws.Cells.Item(291, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(291, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(291, 16).Value = 680
$ws.Cells.Item(291, 17).Value = 25
$ws.Cells.Item(291, 18).Value = "Hortaliza"

# Ensure date column D keeps the date number format used elsewhere in that column.
$ws.Cells.Item(291, 4).NumberFormat = $ws.Cells.Item(292, 4).NumberFormat
